{"js": "// Update the NEP emmeans table: the nutrient-factor levels were reversed\n// (so \"enriched\" is now the reference level) and the model was re-fit,\n// which shifts the lower.CL / upper.CL confidence-interval columns for\n// five of the month rows. Only those ten cells' text content changes;\n// every other cell (month, temperature effect, Q10, SE, df) is untouched.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"values\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected the NEP emmeans table to be present in the document.\");\n}\n\n// row index (as returned by Table.values, 0 = header row) -> { column index: [oldText, newText] }\nconst cellEdits = {\n  1: { 5: [\"0.07\", \"1.07\"], 6: [\"0.11\", \"1.12\"] },  // february\n  2: { 5: [\"-0.02\", \"0.98\"], 6: [\"0.02\", \"1.02\"] },  // june\n  3: { 5: [\"-0.03\", \"0.98\"], 6: [\"0.01\", \"1.01\"] },  // august\n  4: { 5: [\"0.04\", \"1.04\"], 6: [\"0.05\", \"1.06\"] },  // october\n  5: { 5: [\"0.03\", \"1.03\"], 6: [\"0.06\", \"1.07\"] },  // november\n};\n\nfor (const rowIndexStr of Object.keys(cellEdits)) {\n  const rowIndex = Number(rowIndexStr);\n  const columns = cellEdits[rowIndex];\n  for (const colIndexStr of Object.keys(columns)) {\n    const colIndex = Number(colIndexStr);\n    const [oldText, newText] = columns[colIndex];\n\n    const cell = table.getCell(rowIndex, colIndex);\n    // Scope the search to this single cell so identical values elsewhere\n    // in the table (e.g. repeated \"0.01\" SE entries) are left untouched,\n    // then replace in place so the existing run/paragraph formatting\n    // (fonts, size, spacing, keepNext, etc.) is preserved.\n    const found = cell.body.search(oldText, { matchCase: true, matchWholeWord: true });\n    found.load(\"items\");\n    await context.sync();\n\n    if (found.items.length === 0) {\n      throw new Error(`Could not find expected text \"${oldText}\" in row ${rowIndex}, column ${colIndex}.`);\n    }\n\n    found.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the NEP emmeans table: the nutrient-factor levels were reversed\n# (so \"enriched\" is now the reference level) and the model was re-fit,\n# which shifts the lower.CL / upper.CL confidence-interval columns for\n# five of the month rows. Only those ten cells' text content changes;\n# every other cell (month, temperature effect, Q10, SE, df) is untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Table.Cell(row, column) is 1-based and row 1 is the header row, so\n# row 2 = february, row 3 = june, row 4 = august, row 5 = october,\n# row 6 = november. Column 6 = lower.CL, column 7 = upper.CL.\n$edits = @(\n    @{ Row = 2; Col = 6; Old = \"0.07\";  New = \"1.07\" },\n    @{ Row = 2; Col = 7; Old = \"0.11\";  New = \"1.12\" },\n    @{ Row = 3; Col = 6; Old = \"-0.02\"; New = \"0.98\" },\n    @{ Row = 3; Col = 7; Old = \"0.02\";  New = \"1.02\" },\n    @{ Row = 4; Col = 6; Old = \"-0.03\"; New = \"0.98\" },\n    @{ Row = 4; Col = 7; Old = \"0.01\";  New = \"1.01\" },\n    @{ Row = 5; Col = 6; Old = \"0.04\";  New = \"1.04\" },\n    @{ Row = 5; Col = 7; Old = \"0.05\";  New = \"1.06\" },\n    @{ Row = 6; Col = 6; Old = \"0.03\";  New = \"1.03\" },\n    @{ Row = 6; Col = 7; Old = \"0.06\";  New = \"1.07\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $r = $cell.Range\n    # Trim the trailing end-of-cell marker so comparisons/assignment don't\n    # touch it, then verify we are about to overwrite the expected value\n    # before mutating, so the script fails loudly instead of silently\n    # touching the wrong cell if the table layout ever changes.\n    $r.MoveEnd(1, -1) | Out-Null\n    if ($r.Text -ne $edit.Old) {\n        throw \"Row $($edit.Row) Col $($edit.Col): expected '$($edit.Old)' but found '$($r.Text)'\"\n    }\n    $r.Text = $edit.New\n}\n"}
